# "add provider data checked in"
#
# The "Test Results" sheet holds one login's worth of provider rows (rows
# 2+ all reuse the same Username/Password in columns A/B). This edit swaps
# the two existing provider rows for new data and appends a third provider
# as a brand-new row 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 3 down into row 4 first (keeps formatting/type consistent
# with the rest of the table) before the per-cell values are overwritten
# below.
$ws.Range("A3:K3").Copy()
$ws.Range("A4:K4").PasteSpecial()

# --- Row 2 : Helen Kebede ---------------------------------------------
$ws.Range("A2").Value = "srajendran"
$ws.Range("B2").Value = "November@2024!"
$ws.Range("C2").Value = "Helen"
$ws.Range("D2").Value = "Kebede"
$ws.Range("E2").Value = "689 262 7223"
$ws.Range("F2").Value = "W854695"
# "9/30/27" must stay literal text (matches the source column, which is
# stored as text) rather than be auto-parsed into a date serial.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "9/30/27"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = "hkebede@ottersolv.com"
$ws.Range("I2").Value = "Lesly Dorcely"
$ws.Range("J2").Value = "ADJ I"
$ws.Range("K2").Value = "CA0QSVC"

# --- Row 3 : Delonica James --------------------------------------------
$ws.Range("A3").Value = "srajendran"
$ws.Range("B3").Value = "November@2024!"
$ws.Range("C3").Value = "Delonica"
$ws.Range("D3").Value = "James"
$ws.Range("E3").Value = "689 262 7228"
$ws.Range("F3").Value = "G183277"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "djames@ottersolv.com"
$ws.Range("I3").Value = "Hillary Rape"
$ws.Range("J3").Value = "ADJ I"
$ws.Range("K3").Value = "CA0UC2J"

# --- Row 4 (new) : Ashley Hillman ---------------------------------------
$ws.Range("A4").Value = "srajendran"
$ws.Range("B4").Value = "November@2024!"
$ws.Range("C4").Value = "Ashley "
$ws.Range("D4").Value = "Hillman"
$ws.Range("E4").Value = "689 262 7229"
$ws.Range("F4").Value = "W905068"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "ahillman@ottersolv.com"
$ws.Range("I4").Value = "Kristal Fisher"
$ws.Range("J4").Value = "ADJ I"
$ws.Range("K4").Value = "CA0QQH3"

# Grow the sheet's active selection to cover the newly added row.
$ws.Range("K2:K4").Select()
